$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the used range from A1:O15 to A1:P16 and set every cell to 0,
# matching the new assignment-problem matrix for gt/ft matching.
$ws.Range("A1:P16").Value = 0
